# Update Fruta/Hortaliza Chirimoya sheet with corrected weekly data (rows 4-10)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = 44452
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21500
$ws.Range("S4").Value = 2150

# Row 5
$ws.Range("D5").Value = 44448

# Row 6
$ws.Range("D6").Value = 44461
$ws.Range("L6").Value = "Especial"
$ws.Range("N6").Value = 31000
$ws.Range("O6").Value = 32000
$ws.Range("P6").Value = 31500
$ws.Range("S6").Value = 3150

# Row 7
$ws.Range("D7").Value = 44461
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 30000
$ws.Range("O7").Value = 30000
$ws.Range("P7").Value = 30000
$ws.Range("S7").Value = 3000

# Row 8
$ws.Range("D8").Value = 44447

# Row 9
$ws.Range("D9").Value = 44446
$ws.Range("L9").Value = "Primera"
$ws.Range("N9").Value = 21000
$ws.Range("O9").Value = 22000
$ws.Range("P9").Value = 21500
$ws.Range("S9").Value = 2150

# Row 10
$ws.Range("D10").Value = 44487
$ws.Range("N10").Value = 23000
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 23500
$ws.Range("S10").Value = 2350
